$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 43
$ws.Cells.Item($row, 1).Value = "Alessio Farinati"
$ws.Cells.Item($row, 2).Value = "Daniele  Dalbosco | iMontagna"
$ws.Cells.Item($row, 3).Value = "Andrea Conzatti | FC Savignano"
$ws.Cells.Item($row, 4).Value = "ALESSIO FARINATI | Pinguini Trentini"
$ws.Cells.Item($row, 5).Value = "MARCO HEIDEMPERGHER | U.S. Guarna"
$ws.Cells.Item($row, 6).Value = "Federico Rippa | Vigili del Fusto"
